$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where attendance was marked present: Total Attendance Count (D) and Real (E) set to 1
$presentRows = @(4, 10, 11, 13, 17)
foreach ($r in $presentRows) {
    $ws.Cells.Item($r, 4).Value = 1   # D - Total Attendance Count
    $ws.Cells.Item($r, 5).Value = 1   # E - Real
}

# Row 3 is special: Invalid (G) and Absent (H) both set to 1
$ws.Cells.Item(3, 7).Value = 1   # G3 - Invalid
$ws.Cells.Item(3, 8).Value = 1   # H3 - Absent

# Remaining rows: Absent (H) set to 1
$absentRows = @(5, 6, 7, 8, 9, 12, 14, 15, 16, 18)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1   # H - Absent
}
